$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 96, pushing the existing rows 96-100 down to 97-101
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with a new weekly price observation,
# following the same pattern as the surrounding rows in this subset.
$ws.Cells.Item(96, 1).Value = 10
$ws.Cells.Item(96, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(96, 3).Value = "La Araucanía"
$ws.Cells.Item(96, 4).Value = 44516
$ws.Cells.Item(96, 5).Value = 9
$ws.Cells.Item(96, 6).Value = 100112012
$ws.Cells.Item(96, 7).Value = "Espinaca"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 50
$ws.Cells.Item(96, 11).Value = 8000
$ws.Cells.Item(96, 12).Value = 8000
$ws.Cells.Item(96, 13).Value = 8000
$ws.Cells.Item(96, 14).Value = "$/docena de atados"
$ws.Cells.Item(96, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(96, 16).Value = 2667
$ws.Cells.Item(96, 17).Value = 3
$ws.Cells.Item(96, 18).Value = "Hortaliza"
